$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.414.05'
$ws.Range('E2').Value = '  +3.33%  '
$ws.Range('D3').Value = '1.588.17'
$ws.Range('E3').Value = '  +1.17%  '
$ws.Range('E4').Value = '  +1.04%  '
$ws.Range('D5').Value = '213.45'
$ws.Range('E5').Value = '  +1.15%  '
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('E7').Value = '  +1.10%  '
$ws.Range('D8').Value = '24.20'
$ws.Range('E8').Value = '  +6.87%  '
$ws.Range('E9').Value = '  +0.17%  '
$ws.Range('D10').Value = '0.0600'
$ws.Range('E10').Value = '  +0.54%  '
$ws.Range('E11').Value = '  +1.87%  '
$ws.Range('D12').Value = '1.815.90'
$ws.Range('E12').Value = '  +1.43%  '
$ws.Range('D13').Value = '1.587.63'
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('E14').Value = '  +1.93%  '
$ws.Range('E15').Value = '  -0.67%  '
$ws.Range('D16').Value = '28.422.34'
$ws.Range('E16').Value = '  +3.35%  '
$ws.Range('D17').Value = '63.18'
$ws.Range('E17').Value = '  +1.15%  '
$ws.Range('D18').Value = '229.41'
$ws.Range('E18').Value = '  +1.28%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0706'
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = '7.48'
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  +1.07%  '
$ws.Range('E22').Value = '  -1.61%  '
$ws.Range('D23').Value = '9.33'
$ws.Range('E23').Value = '  -0.94%  '
$ws.Range('E24').Value = '  +0.19%  '
$ws.Range('D25').Value = '151.80'
$ws.Range('E25').Value = '  +1.54%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E27').Value = '  -0.84%  '
$ws.Range('E28').Value = '  -0.86%  '
$ws.Range('E29').Value = '  +1.06%  '
$ws.Range('D30').Value = '1.13'
$ws.Range('E30').Value = '  -0.53%  '
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').Value = '3.25'
$ws.Range('E32').Value = '  +0.18%  '
$ws.Range('E33').Value = '  +0.53%  '
$ws.Range('D34').Value = '1.399.33'
$ws.Range('E34').Value = '  -3.63%  '
$ws.Range('E35').Value = '  -1.19%  '
$ws.Range('E36').Value = '  -8.78%  '
$ws.Range('E37').Value = '  +1.54%  '
$ws.Range('D38').Value = '0.0166'
$ws.Range('E38').Value = '  -0.78%  '
$ws.Range('E39').Value = '  +8.73%  '
$ws.Range('E40').Value = '  -0.43%  '
$ws.Range('D41').Value = '0.810'
$ws.Range('E41').Value = '  -0.72%  '
$ws.Range('E42').Value = '  +1.08%  '
$ws.Range('E43').Value = '  +2.54%  '
$ws.Range('E44').Value = '  -2.96%  '
$ws.Range('E45').Value = '  +0.67%  '
$ws.Range('E46').Value = '  -0.67%  '
$ws.Range('D47').Value = '1.725.92'
$ws.Range('E47').Value = '  +1.42%  '
$ws.Range('B48').Value = 'mCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D48').Value = '2.14'
$ws.Range('E48').Value = '  +2.30%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '87.19'
$ws.Range('E49').Value = '  +0.56%  '
$ws.Range('D50').Value = '0.0₆0104'
$ws.Range('E50').Value = '  +2.75%  '
$ws.Range('D51').Value = '0.0519'
$ws.Range('E51').Value = '  -1.06%  '
